$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, shifting existing rows 150-251 down to 151-252
$ws.Rows("150:150").Insert()

# Populate the newly inserted row 150 with the new record
$ws.Range("A150").Value = 4
$ws.Range("B150").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C150").Value = "Los Lagos"
$ws.Range("D150").Value = 44827
$ws.Range("E150").Value = 10
$ws.Range("F150").Value = "Fruta"
$ws.Range("G150").Value = 100108
$ws.Range("H150").Value = "Tropicales y subtropicales"
$ws.Range("I150").Value = 100108002
$ws.Range("J150").Value = "Mango"
$ws.Range("K150").Value = "Sin especificar"
$ws.Range("L150").Value = "Primera"
$ws.Range("M150").Value = 120
$ws.Range("N150").Value = 9000
$ws.Range("O150").Value = 10000
$ws.Range("P150").Value = 9500
$ws.Range("Q150").Value = "$/bandeja 4 kilos"
$ws.Range("R150").Value = "Brasil"
$ws.Range("S150").Value = 2375
$ws.Range("T150").Value = 4
